# citgo_locations_by_state.xlsx — "Add files via upload" edit
#
# Target state (per the OOXML diff):
#   - Sheet "citgo"   -> stays visible / active, selection moves to G15
#   - Sheet "Sheet1"  -> becomes hidden, selection (kept for when shown) moves to C14
#   - Sheet "Abbreviations" -> becomes hidden (selection left untouched)
#   - bookViews/workbookView gains firstSheet="1" (tab-strip scroll position;
#     purely cosmetic UI state, not reproducible through this object model)

$wb = $excel.ActiveWorkbook

# --- Update the remembered selection on "Sheet1" (A14 -> C14) ---
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Activate()
$sheet1.Range("C14").Select() | Out-Null

# --- Update the remembered selection on "citgo" (A1:C31/C31 -> G15) ---
$citgo = $wb.Worksheets.Item("citgo")
$citgo.Activate()
$citgo.Range("G15").Select() | Out-Null

# --- Hide the helper sheets, leaving "citgo" as the sole visible/active tab ---
$sheet1.Visible = $false
$wb.Worksheets.Item("Abbreviations").Visible = $false
